# Add a new person row to the lab list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 is intentionally left blank; the new entry goes in row 15, column A only.
$ws.Range("A15").Value = "Here is an extra person"

# Match the author's final selection (cell A15) as recorded in the saved file.
$ws.Range("A15").Select() | Out-Null
